# Apply the commit "Added periodic & upfront related scenarios":
#  - ProductLoanInput!B17 ("repaymentstrategy" value) changes from "Mifos style"
#    to "Penalties, Fees, Interest, Principal order", with new left/top aligned
#    formatting.
#  - The ProductLoanInput sheet becomes the active/selected tab (with B17
#    selected), replacing ProductLoanOutput as the active tab.

$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")

$cell = $wsInput.Range("B17")
$cell.Value = "Penalties, Fees, Interest, Principal order"
$cell.HorizontalAlignment = -4131
$cell.VerticalAlignment = -4160

$wsInput.Activate()
$wsInput.Range("B17").Select() | Out-Null
